$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sorted language/value data (rows 2-21), descending by value
$data = @(
    @("English", 29.46028061126449),
    @("Spanish", 9.367558929096459),
    @("Japanese", 8.474824939210009),
    @("German", 7.591510362179372),
    @("Arabic", 5.454635801172345),
    @("French", 4.643752854412187),
    @("Italian", 4.445903922931991),
    @("Portuguese", 4.316655195967407),
    @("Chinese", 3.578221562981212),
    @("Persian", 2.075595326933791),
    @("Dutch", 1.948150317324575),
    @("Malay-Indonesian", 1.945403479783394),
    @("Turkish", 1.281499919481114),
    @("Polish", 1.101971097610346),
    @("Korean", 0.7756102518107317),
    @("Urdu", 0.6859195221398643),
    @("Swedish", 0.6746904974275217),
    @("Thai", 0.6262318033801899),
    @("Bengali", 0.4092139541233332),
    @("Vietnamese", 0.2543595141139115)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the old trailing rows (Russian, Uzbek at rows 22 and 23)
$ws.Range("A22:B23").Delete()
